# Power_Calc: BOM updated - one more component added in the second supply
# rail's item list (O8: 4 -> 5 units), selection moved to O8 to reflect
# where the edit was made.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O8").Value = 5
$ws.Range("O8").Select()
